$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the dataset. It belongs right
# before the existing row 205, so insert a blank row there; Excel shifts
# the old rows 205:237 down to 206:238 (and grows the used range to R238).
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205 with the new record's data.
$ws.Cells.Item(205, 1).Value = 9
$ws.Cells.Item(205, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(205, 3).Value = "Metropolitana"
$ws.Cells.Item(205, 4).Value = 44522
$ws.Cells.Item(205, 5).Value = 13
$ws.Cells.Item(205, 6).Value = 100112052
$ws.Cells.Item(205, 7).Value = "Albahaca"
$ws.Cells.Item(205, 8).Value = "Sin especificar"
$ws.Cells.Item(205, 9).Value = "Primera"
$ws.Cells.Item(205, 10).Value = 160
$ws.Cells.Item(205, 11).Value = 4000
$ws.Cells.Item(205, 12).Value = 5000
$ws.Cells.Item(205, 13).Value = 4500
$ws.Cells.Item(205, 14).Value = "$/paquete"
$ws.Cells.Item(205, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(205, 16).Value = 4500
$ws.Cells.Item(205, 17).Value = 1
$ws.Cells.Item(205, 18).Value = "Hortaliza"
